$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new components at the end of the list (rows 18 and 19)
$ws.Range("A18").Value = 18
$ws.Range("B18").Value = "метилциклопентан"
$ws.Range("A19").Value = 19
$ws.Range("B19").Value = "циклогексан"

# Capitalize the first three component names (пропан, метан, этан -> Пропан, Метан, Этан)
$ws.Range("B3").Value = "Пропан"
$ws.Range("B1").Value = "Метан"
$ws.Range("B2").Value = "Этан"

# Update the selected cell to match the saved workbook state
$ws.Range("B3").Select()
